# Apply the content/formatting edits described in the commit:
#  - correct a few contact-record values (names, email typo fix, trim
#    trailing commas from the hobbies list)
#  - make the Pincode/Phone data font color explicit black
#  - bump the header/data row height slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (first contact) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# --- Row 3 (second contact) ---
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Explicit black font for the Pincode/Phone number cells ---
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0

# --- Row height bump for the header + two data rows ---
$ws.Rows("1:3").RowHeight = 19.5
